# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E, rows 16-50) is re-ordered from descending
# (2003 .. 1705) to ascending (1705 .. 2003), and the corresponding
# "Valor Mora" column (F, rows 16-50) is updated so periods through 1811
# carry 29509 and periods from 1812 onward carry 31249.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=16; E="1705"; F=29509},
    @{Row=17; E="1706"; F=29509},
    @{Row=18; E="1707"; F=29509},
    @{Row=19; E="1708"; F=29509},
    @{Row=20; E="1709"; F=29509},
    @{Row=21; E="1710"; F=29509},
    @{Row=22; E="1711"; F=29509},
    @{Row=23; E="1712"; F=29509},
    @{Row=24; E="1801"; F=29509},
    @{Row=25; E="1802"; F=29509},
    @{Row=26; E="1803"; F=29509},
    @{Row=27; E="1804"; F=29509},
    @{Row=28; E="1805"; F=29509},
    @{Row=29; E="1806"; F=29509},
    @{Row=30; E="1807"; F=29509},
    @{Row=31; E="1808"; F=29509},
    @{Row=32; E="1809"; F=29509},
    @{Row=33; E="1810"; F=29509},
    @{Row=34; E="1811"; F=29509},
    @{Row=35; E="1812"; F=31249},
    @{Row=36; E="1901"; F=31249},
    @{Row=37; E="1902"; F=31249},
    @{Row=38; E="1903"; F=31249},
    @{Row=39; E="1904"; F=31249},
    @{Row=40; E="1905"; F=31249},
    @{Row=41; E="1906"; F=31249},
    @{Row=42; E="1907"; F=31249},
    @{Row=43; E="1908"; F=31249},
    @{Row=44; E="1909"; F=31249},
    @{Row=45; E="1910"; F=31249},
    @{Row=46; E="1911"; F=31249},
    @{Row=47; E="1912"; F=31249},
    @{Row=48; E="2001"; F=31249},
    @{Row=49; E="2002"; F=31249},
    @{Row=50; E="2003"; F=31249}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 6).Value = $u.F
}
